# Colours.xlsx — rename the "UNC in coincard" row to reflect that the
# colour-coding example also covers coins in other forms of original sealing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# D7 held the example text "UNC in coincard" for the Yellow colour row;
# broaden it to also cover other original sealing.
$ws.Range("D7").Value = "UNC in coincard or in other original sealing"

# Restore the cursor/selection position as last left by the author.
$ws.Range("D29").Select()
